$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three newly scraped part numbers to the bottom of column A
# (sharedStrings gains RXG11RD, RXM2AB1B7, RXM4AB1B7; dimension grows to A1:A20)
$ws.Range("A18").Value = "RXG11RD"
$ws.Range("A19").Value = "RXM2AB1B7"
$ws.Range("A20").Value = "RXM4AB1B7"

# Widen column A to fit the longer part numbers (stored width ends up 12.5)
$ws.Columns("A").ColumnWidth = 11.6

# Move the selection down past the newly added data, matching the
# author's post-edit cursor position (A21 active cell, A21:B56 selected)
$ws.Range("A21:B56").Select() | Out-Null
